$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("G2:G5").Value = "2016-10-26 08:35:22"
$ov.Range("E3:F3").Value = "In Translation"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("E2:E5").Value = "mt"
$zh.Range("H2:H5").Value = "2016-10-26 08:35:10"
$zh.Range("C3").Value = "In Translation"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("E2:E5").Value = "mt"
$de.Range("H2:H5").Value = "2016-10-26 08:35:22"
$de.Range("C3").Value = "In Translation"
